# Refresh the crypto price/volume table to the latest scraped snapshot.
# Price cells (column D) sometimes hold purely numeric-looking text (e.g. "1.013",
# "88.77") that must stay literal text (matching the source feed's inline-string
# cells) instead of being auto-coerced into a floating point number by Excel's
# normal cell-entry parsing. For those cells we force the "@" (Text) number
# format before writing the value so the exact original text is preserved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bitcoin
$ws.Range('D2').Value = '26.980.97'
$ws.Range('E2').Value = '  +1.28%  '
# Ethereum
$ws.Range('D3').Value = '1.847.30'
$ws.Range('E3').Value = '  +1.14%  '
# TetherUSD
$ws.Range('E4').Value = '  +0.43%  '
# row 5 (was USDC, now BNB)
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.10'
$ws.Range('E5').Value = '  +0.15%  '
# row 6 (was BNB, now USDC)
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.012'
$ws.Range('E6').Value = '  +0.36%  '
# XRP
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4781'
$ws.Range('E7').Value = '  +2.51%  '
# Cardano
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3682'
$ws.Range('E8').Value = '  +2.31%  '
# Dogecoin
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07231'
$ws.Range('E9').Value = '  +1.15%  '
# Polygon
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9300'
$ws.Range('E10').Value = '  +2.81%  '
# Solana
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.73'
$ws.Range('E11').Value = '  +1.50%  '
# TRON
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07716'
$ws.Range('E12').Value = '  +0.19%  '
# WrappedEther
$ws.Range('D13').Value = '1.789.50'
$ws.Range('E13').Value = '  -2.30%  '
# Polkadot
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.336'
$ws.Range('E14').Value = '  +1.22%  '
# Chainlink
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.432'
$ws.Range('E15').Value = '  +0.92%  '
# Litecoin
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.77'
$ws.Range('E16').Value = '  +1.24%  '
# BinanceUSD
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.015'
$ws.Range('E17').Value = '  +0.48%  '
# ShibaInu
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008637'
$ws.Range('E18').Value = '  +0.86%  '
# Dai
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.011'
$ws.Range('E19').Value = '  +0.37%  '
# WrappedBTC
$ws.Range('D20').Value = '27.021.99'
$ws.Range('E20').Value = '  +1.29%  '
# Avalanche
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.49'
$ws.Range('E21').Value = '  +1.90%  '
# Uniswap
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.056'
$ws.Range('E22').Value = '  +0.52%  '
# Cosmos
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.64'
$ws.Range('E23').Value = '  +0.76%  '
# Toncoin
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.929'
$ws.Range('E24').Value = '  +1.19%  '
# Monero
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.75'
# EthereumClassic
$ws.Range('E26').Value = '  +1.42%  '
# LidoDAOToken
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.999'
$ws.Range('E27').Value = '  +0.64%  '
# BitcoinCash
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.36'
$ws.Range('E28').Value = '  +0.34%  '
# InternetComputer(DFINITY)
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.983'
$ws.Range('E29').Value = '  +2.31%  '
# Stellar
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08900'
$ws.Range('E30').Value = '  +0.96%  '
# HuobiToken
$ws.Range('E31').Value = '  +5.36%  '
# ARBITRUM
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.175'
$ws.Range('E32').Value = '  +0.62%  '
# ImmutableX
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7434'
$ws.Range('E33').Value = '  +1.16%  '
# Filecoin
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.504'
$ws.Range('E34').Value = '  +1.42%  '
# RenderToken
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.739'
$ws.Range('E35').Value = '  -3.23%  '
# TrustWalletToken
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.118'
$ws.Range('E36').Value = '  +3.44%  '
# VeChain
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01958'
$ws.Range('E37').Value = '  +1.22%  '
# Hedera
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05263'
$ws.Range('E38').Value = '  +1.93%  '
# MXToken
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.977'
$ws.Range('E39').Value = '  +1.93%  '
# TheSandbox
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5220'
$ws.Range('E40').Value = '  +2.96%  '
# FraxShare
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.995'
$ws.Range('E41').Value = '  +1.70%  '
# Algorand
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1514'
$ws.Range('E42').Value = '  +1.05%  '
# Aptos
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.204'
$ws.Range('E43').Value = '  +1.57%  '
# EnergySwap
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.65'
$ws.Range('E44').Value = '  +6.34%  '
# Decentraland
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4748'
$ws.Range('E45').Value = '  +1.67%  '
# PaxDollar
$ws.Range('E46').Value = '  +0.42%  '
# Quant
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.03'
$ws.Range('E47').Value = '  +3.84%  '
# NEARProtocol
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.610'
$ws.Range('E48').Value = '  +2.27%  '
# Aave
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '66.05'
$ws.Range('E49').Value = '  +3.18%  '
# Cronos
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06059'
$ws.Range('E50').Value = '  +0.39%  '
# EOS
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8870'
$ws.Range('E51').Value = '  +3.88%  '
